$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# --- Set guaranteed dispatch values (B column) for the fuel sources that ---
# --- now carry a value of 1, and fill the formulas across the year columns ---
$rowsToOne = 4,5,6,7,8,9,10
foreach ($r in $rowsToOne) {
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Formula = "=`$B$r"
    $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 37)).Formula = "=`$B$r"
}

# --- Add new rows 13-17 for additional fuel sources ---

# Row 13: lignite (all zero, formulas mirror the $B<row> pattern)
$ws.Range("A13").Value = "lignite"
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Formula = "=`$B13"
$ws.Range($ws.Cells.Item(13, 4), $ws.Cells.Item(13, 37)).Formula = "=`$B13"

# --- Rename existing fuel source labels ---
$ws.Range("A2").Value = "hard coal"
$ws.Range("A6").Value = "onshore wind"

# Row 14: offshore wind (value 1, formulas mirror the $B<row> pattern)
$ws.Range("A14").Value = "offshore wind"
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Formula = "=`$B14"
$ws.Range($ws.Cells.Item(14, 4), $ws.Cells.Item(14, 37)).Formula = "=`$B14"

# Row 15: crude oil (plain static zeros, no formulas)
$ws.Range("A15").Value = "crude oil"
$ws.Range($ws.Cells.Item(15, 2), $ws.Cells.Item(15, 37)).Value = 0

# Row 16: heavy or residual fuel oil (plain static zeros, no formulas)
$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range($ws.Cells.Item(16, 2), $ws.Cells.Item(16, 37)).Value = 0

# Row 17: municipal solid waste (value 1, chained formulas referencing the previous column)
$ws.Range("A17").Value = "municipal solid waste"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Formula = "=B17"
$ws.Range($ws.Cells.Item(17, 4), $ws.Cells.Item(17, 37)).Formula = "=C17"

# --- Header cell (A1): new label, bold + wrap text, taller row ---
$ws.Range("A1").Value = "BAU Guaranteed Dispatch (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

$wb.Save()
